$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the two narrative rows that were dropped from the use-case
#     ("3. Seleciona personalizar especificações" / "4. Mostra opções").
#     This shifts everything below up by two rows and keeps merged
#     cells / row styles intact, matching rows 11-27 -> 9-25.
$ws.Rows("9:10").Delete()

# --- Pré condição / Pós condição text (row 4 / row 5)
$ws.Range("C4").Value = "Autenticado no sistema e escolheu personalizar o carro"
$ws.Range("C5").Value = "Carro encomendado"

# --- Renumbered narrative steps (rows 7-16 after the deletion above)
$ws.Range("C7").Value = "1.<<include>> Comprar Carro"
$ws.Range("D8").Value = "2. Mostra opções"
$ws.Range("C9").Value = "3. Escolhe especificações"
$ws.Range("D10").Value = "4. Regista especificações"
$ws.Range("D11").Value = "5. Verifica especificações"
$ws.Range("D12").Value = "6. Calcula preço"
$ws.Range("D13").Value = "7. Mostra preço"
$ws.Range("D14").Value = "8. Pergunta se pretende confirmar"
$ws.Range("C15").Value = "9. Confirma"
$ws.Range("D16").Value = "10. Regista no sistema e adiciona a fila de espera"

# --- Alternativa 1 block (rows 18-20)
$ws.Range("B18").Value = "Alternativa 1 [Peças Incompativeis] (passo 5)"
$ws.Range("D18").Value = "5.1 Verifica que especificações não estão corretas"
$ws.Range("D19").Value = "5.2 Informa que escolheu peças incompativeis e/ou peças em falta"
$ws.Range("D20").Value = "Regressa a 3"

# --- Exceção 1 block (rows 21-23): renamed from "Alternativa 2", and the
#     trailing "Regressa a 1" note in D22 is removed.
$ws.Range("B21").Value = "Exceção 1 [Não confirma a compra] (passo 9)"
$ws.Range("C21").Value = "9.1 Não confirma "
$ws.Range("D22").Value = ""

# --- Sheet view: select C5:D5 (also clears the old topLeftCell scroll state)
$ws.Range("C5:D5").Select()
